$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text / content changes (commit: "VLJ # in spreadsheet instead of css_id") ---

# Header: "CSS Id" -> "VLJ #"
$ws.Range("C2").Value = "VLJ #"

# Example / Bernard Jones rows: id "BVAJONESB" -> "123"
$ws.Range("C3").Value = "123"
$ws.Range("C4").Value = "123"
$ws.Range("C5").Value = "123"
$ws.Range("C6").Value = "123"
$ws.Range("C7").Value = "123"

# Second judge: "Roth, Lauren" / "DSUSER" -> "Huels, Stuart" / "860"
$ws.Range("B8").Value = "Huels, Stuart"
$ws.Range("C8").Value = "860"
$ws.Range("B9").Value = "Huels, Stuart"
$ws.Range("C9").Value = "860"

# --- Structural change: append a new (blank) trailing row 10, matching the ---
# --- look of the preceding row, extending the table/dimension to A1:I10.  ---
$ws.Range("A9:I9").Copy()
$ws.Range("A10:I10").PasteSpecial(-4122)
$ws.Range("A10:I10").ClearContents()
$excel.CutCopyMode = 0
$ws.Range("A10").RowHeight = 17
